$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 25
$ws.Range("F3").Value = 8
$ws.Range("F4").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("F6").Value = 0
$ws.Range("F7").Value = 0
$ws.Range("F8").Value = 0
$ws.Range("F9").Value = 26
$ws.Range("F10").Value = 9
$ws.Range("F11").Value = 1
$ws.Range("F12").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("F14").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("F16").Value = 27
$ws.Range("F17").Value = 11
$ws.Range("F18").Value = 1
$ws.Range("F19").Value = 0
$ws.Range("F20").Value = 0
$ws.Range("F21").Value = 0
$ws.Range("F22").Value = 0
$ws.Range("F23").Value = 27
$ws.Range("F24").Value = 10
$ws.Range("F25").Value = 0
$ws.Range("F26").Value = 0
$ws.Range("F27").Value = 0
$ws.Range("F28").Value = 0
$ws.Range("F29").Value = 0
$ws.Range("F30").Value = 41
$ws.Range("F31").Value = 29
$ws.Range("F32").Value = 9
$ws.Range("F33").Value = 8
$ws.Range("F34").Value = 5
$ws.Range("F35").Value = 1
$ws.Range("F36").Value = 0
$ws.Range("F37").Value = 41
$ws.Range("F38").Value = 25
$ws.Range("F39").Value = 8
$ws.Range("F40").Value = 7
$ws.Range("F41").Value = 3
$ws.Range("F42").Value = 2
$ws.Range("F43").Value = 0
$ws.Range("F44").Value = 42
$ws.Range("F45").Value = 29
$ws.Range("F46").Value = 10
$ws.Range("F47").Value = 7
$ws.Range("F48").Value = 4
$ws.Range("F49").Value = 2
$ws.Range("F50").Value = 0
$ws.Range("F51").Value = 42
$ws.Range("F52").Value = 28
$ws.Range("F53").Value = 9
$ws.Range("F54").Value = 5
$ws.Range("F55").Value = 3
$ws.Range("F56").Value = 1
$ws.Range("F57").Value = 0
